# The "Oo" / "Oa" / "Of" rows (and their matching columns) are being
# removed from the dependency matrix. Deleting the whole row/column shifts
# everything below/right of them up and to the left, which is exactly what
# the target file shows (table shrinks from A1:T20 to A1:Q17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Oo" lives in column H (header H1) / row 8 (label A8)
# "Oa" lives in column I (header I1) / row 9 (label A9)
# "Of" lives in column J (header J1) / row 10 (label A10)
$ws.Range("H1:J1").EntireColumn.Delete()
$ws.Range("A8:A10").EntireRow.Delete()

# Update the selection to match the post-edit view.
[void]$ws.Range("K21").Select()
